$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 161 <-> row 163 for columns B,F,G
$tmp = $ws.Range("B161").Value2
$ws.Range("B161").Value = $ws.Range("B163").Value2
$ws.Range("B163").Value = $tmp
$tmp = $ws.Range("F161").Value2
$ws.Range("F161").Value = $ws.Range("F163").Value2
$ws.Range("F163").Value = $tmp
$tmp = $ws.Range("G161").Value2
$ws.Range("G161").Value = $ws.Range("G163").Value2
$ws.Range("G163").Value = $tmp

# Swap row 183 <-> row 184 for columns B,E,F,G
$tmp = $ws.Range("B183").Value2
$ws.Range("B183").Value = $ws.Range("B184").Value2
$ws.Range("B184").Value = $tmp
$tmp = $ws.Range("E183").Value2
$ws.Range("E183").Value = $ws.Range("E184").Value2
$ws.Range("E184").Value = $tmp
$tmp = $ws.Range("F183").Value2
$ws.Range("F183").Value = $ws.Range("F184").Value2
$ws.Range("F184").Value = $tmp
$tmp = $ws.Range("G183").Value2
$ws.Range("G183").Value = $ws.Range("G184").Value2
$ws.Range("G184").Value = $tmp

# Swap row 279 <-> row 280 for columns B,E,F,G
$tmp = $ws.Range("B279").Value2
$ws.Range("B279").Value = $ws.Range("B280").Value2
$ws.Range("B280").Value = $tmp
$tmp = $ws.Range("E279").Value2
$ws.Range("E279").Value = $ws.Range("E280").Value2
$ws.Range("E280").Value = $tmp
$tmp = $ws.Range("F279").Value2
$ws.Range("F279").Value = $ws.Range("F280").Value2
$ws.Range("F280").Value = $tmp
$tmp = $ws.Range("G279").Value2
$ws.Range("G279").Value = $ws.Range("G280").Value2
$ws.Range("G280").Value = $tmp

# Swap row 313 <-> row 314 for columns B,F,G
$tmp = $ws.Range("B313").Value2
$ws.Range("B313").Value = $ws.Range("B314").Value2
$ws.Range("B314").Value = $tmp
$tmp = $ws.Range("F313").Value2
$ws.Range("F313").Value = $ws.Range("F314").Value2
$ws.Range("F314").Value = $tmp
$tmp = $ws.Range("G313").Value2
$ws.Range("G313").Value = $ws.Range("G314").Value2
$ws.Range("G314").Value = $tmp

# Swap row 316 <-> row 317 for columns B,D,E,F,G
$tmp = $ws.Range("B316").Value2
$ws.Range("B316").Value = $ws.Range("B317").Value2
$ws.Range("B317").Value = $tmp
$tmp = $ws.Range("D316").Value2
$ws.Range("D316").Value = $ws.Range("D317").Value2
$ws.Range("D317").Value = $tmp
$tmp = $ws.Range("E316").Value2
$ws.Range("E316").Value = $ws.Range("E317").Value2
$ws.Range("E317").Value = $tmp
$tmp = $ws.Range("F316").Value2
$ws.Range("F316").Value = $ws.Range("F317").Value2
$ws.Range("F317").Value = $tmp
$tmp = $ws.Range("G316").Value2
$ws.Range("G316").Value = $ws.Range("G317").Value2
$ws.Range("G317").Value = $tmp

# Swap row 351 <-> row 352 for columns B,E,F,G
$tmp = $ws.Range("B351").Value2
$ws.Range("B351").Value = $ws.Range("B352").Value2
$ws.Range("B352").Value = $tmp
$tmp = $ws.Range("E351").Value2
$ws.Range("E351").Value = $ws.Range("E352").Value2
$ws.Range("E352").Value = $tmp
$tmp = $ws.Range("F351").Value2
$ws.Range("F351").Value = $ws.Range("F352").Value2
$ws.Range("F352").Value = $tmp
$tmp = $ws.Range("G351").Value2
$ws.Range("G351").Value = $ws.Range("G352").Value2
$ws.Range("G352").Value = $tmp

# Swap row 372 <-> row 373 for columns B,E,F,G
$tmp = $ws.Range("B372").Value2
$ws.Range("B372").Value = $ws.Range("B373").Value2
$ws.Range("B373").Value = $tmp
$tmp = $ws.Range("E372").Value2
$ws.Range("E372").Value = $ws.Range("E373").Value2
$ws.Range("E373").Value = $tmp
$tmp = $ws.Range("F372").Value2
$ws.Range("F372").Value = $ws.Range("F373").Value2
$ws.Range("F373").Value = $tmp
$tmp = $ws.Range("G372").Value2
$ws.Range("G372").Value = $ws.Range("G373").Value2
$ws.Range("G373").Value = $tmp

# Swap row 379 <-> row 380 for columns B,E,F,G
$tmp = $ws.Range("B379").Value2
$ws.Range("B379").Value = $ws.Range("B380").Value2
$ws.Range("B380").Value = $tmp
$tmp = $ws.Range("E379").Value2
$ws.Range("E379").Value = $ws.Range("E380").Value2
$ws.Range("E380").Value = $tmp
$tmp = $ws.Range("F379").Value2
$ws.Range("F379").Value = $ws.Range("F380").Value2
$ws.Range("F380").Value = $tmp
$tmp = $ws.Range("G379").Value2
$ws.Range("G379").Value = $ws.Range("G380").Value2
$ws.Range("G380").Value = $tmp

# Swap row 382 <-> row 383 for columns B,E,F,G
$tmp = $ws.Range("B382").Value2
$ws.Range("B382").Value = $ws.Range("B383").Value2
$ws.Range("B383").Value = $tmp
$tmp = $ws.Range("E382").Value2
$ws.Range("E382").Value = $ws.Range("E383").Value2
$ws.Range("E383").Value = $tmp
$tmp = $ws.Range("F382").Value2
$ws.Range("F382").Value = $ws.Range("F383").Value2
$ws.Range("F383").Value = $tmp
$tmp = $ws.Range("G382").Value2
$ws.Range("G382").Value = $ws.Range("G383").Value2
$ws.Range("G383").Value = $tmp

# Swap row 389 <-> row 390 for columns B,F,G
$tmp = $ws.Range("B389").Value2
$ws.Range("B389").Value = $ws.Range("B390").Value2
$ws.Range("B390").Value = $tmp
$tmp = $ws.Range("F389").Value2
$ws.Range("F389").Value = $ws.Range("F390").Value2
$ws.Range("F390").Value = $tmp
$tmp = $ws.Range("G389").Value2
$ws.Range("G389").Value = $ws.Range("G390").Value2
$ws.Range("G390").Value = $tmp

# Swap row 431 <-> row 432 for columns B,C,F,G
$tmp = $ws.Range("B431").Value2
$ws.Range("B431").Value = $ws.Range("B432").Value2
$ws.Range("B432").Value = $tmp
$tmp = $ws.Range("C431").Value2
$ws.Range("C431").Value = $ws.Range("C432").Value2
$ws.Range("C432").Value = $tmp
$tmp = $ws.Range("F431").Value2
$ws.Range("F431").Value = $ws.Range("F432").Value2
$ws.Range("F432").Value = $tmp
$tmp = $ws.Range("G431").Value2
$ws.Range("G431").Value = $ws.Range("G432").Value2
$ws.Range("G432").Value = $tmp

# Swap row 457 <-> row 458 for columns B,E,F,G
$tmp = $ws.Range("B457").Value2
$ws.Range("B457").Value = $ws.Range("B458").Value2
$ws.Range("B458").Value = $tmp
$tmp = $ws.Range("E457").Value2
$ws.Range("E457").Value = $ws.Range("E458").Value2
$ws.Range("E458").Value = $tmp
$tmp = $ws.Range("F457").Value2
$ws.Range("F457").Value = $ws.Range("F458").Value2
$ws.Range("F458").Value = $tmp
$tmp = $ws.Range("G457").Value2
$ws.Range("G457").Value = $ws.Range("G458").Value2
$ws.Range("G458").Value = $tmp

# Swap row 581 <-> row 582 for columns B,E,F,G
$tmp = $ws.Range("B581").Value2
$ws.Range("B581").Value = $ws.Range("B582").Value2
$ws.Range("B582").Value = $tmp
$tmp = $ws.Range("E581").Value2
$ws.Range("E581").Value = $ws.Range("E582").Value2
$ws.Range("E582").Value = $tmp
$tmp = $ws.Range("F581").Value2
$ws.Range("F581").Value = $ws.Range("F582").Value2
$ws.Range("F582").Value = $tmp
$tmp = $ws.Range("G581").Value2
$ws.Range("G581").Value = $ws.Range("G582").Value2
$ws.Range("G582").Value = $tmp

# Swap row 583 <-> row 584 for columns B,E,F,G
$tmp = $ws.Range("B583").Value2
$ws.Range("B583").Value = $ws.Range("B584").Value2
$ws.Range("B584").Value = $tmp
$tmp = $ws.Range("E583").Value2
$ws.Range("E583").Value = $ws.Range("E584").Value2
$ws.Range("E584").Value = $tmp
$tmp = $ws.Range("F583").Value2
$ws.Range("F583").Value = $ws.Range("F584").Value2
$ws.Range("F584").Value = $tmp
$tmp = $ws.Range("G583").Value2
$ws.Range("G583").Value = $ws.Range("G584").Value2
$ws.Range("G584").Value = $tmp

# Swap row 586 <-> row 587 for columns B,E,F,G
$tmp = $ws.Range("B586").Value2
$ws.Range("B586").Value = $ws.Range("B587").Value2
$ws.Range("B587").Value = $tmp
$tmp = $ws.Range("E586").Value2
$ws.Range("E586").Value = $ws.Range("E587").Value2
$ws.Range("E587").Value = $tmp
$tmp = $ws.Range("F586").Value2
$ws.Range("F586").Value = $ws.Range("F587").Value2
$ws.Range("F587").Value = $tmp
$tmp = $ws.Range("G586").Value2
$ws.Range("G586").Value = $ws.Range("G587").Value2
$ws.Range("G587").Value = $tmp

# Swap row 593 <-> row 594 for columns B,E,F,G
$tmp = $ws.Range("B593").Value2
$ws.Range("B593").Value = $ws.Range("B594").Value2
$ws.Range("B594").Value = $tmp
$tmp = $ws.Range("E593").Value2
$ws.Range("E593").Value = $ws.Range("E594").Value2
$ws.Range("E594").Value = $tmp
$tmp = $ws.Range("F593").Value2
$ws.Range("F593").Value = $ws.Range("F594").Value2
$ws.Range("F594").Value = $tmp
$tmp = $ws.Range("G593").Value2
$ws.Range("G593").Value = $ws.Range("G594").Value2
$ws.Range("G594").Value = $tmp

# Swap row 601 <-> row 602 for columns B,E,F,G
$tmp = $ws.Range("B601").Value2
$ws.Range("B601").Value = $ws.Range("B602").Value2
$ws.Range("B602").Value = $tmp
$tmp = $ws.Range("E601").Value2
$ws.Range("E601").Value = $ws.Range("E602").Value2
$ws.Range("E602").Value = $tmp
$tmp = $ws.Range("F601").Value2
$ws.Range("F601").Value = $ws.Range("F602").Value2
$ws.Range("F602").Value = $tmp
$tmp = $ws.Range("G601").Value2
$ws.Range("G601").Value = $ws.Range("G602").Value2
$ws.Range("G602").Value = $tmp

# Swap row 604 <-> row 605 for columns B,E,F,G
$tmp = $ws.Range("B604").Value2
$ws.Range("B604").Value = $ws.Range("B605").Value2
$ws.Range("B605").Value = $tmp
$tmp = $ws.Range("E604").Value2
$ws.Range("E604").Value = $ws.Range("E605").Value2
$ws.Range("E605").Value = $tmp
$tmp = $ws.Range("F604").Value2
$ws.Range("F604").Value = $ws.Range("F605").Value2
$ws.Range("F605").Value = $tmp
$tmp = $ws.Range("G604").Value2
$ws.Range("G604").Value = $ws.Range("G605").Value2
$ws.Range("G605").Value = $tmp

# Swap row 715 <-> row 716 for columns B,E,F,G
$tmp = $ws.Range("B715").Value2
$ws.Range("B715").Value = $ws.Range("B716").Value2
$ws.Range("B716").Value = $tmp
$tmp = $ws.Range("E715").Value2
$ws.Range("E715").Value = $ws.Range("E716").Value2
$ws.Range("E716").Value = $tmp
$tmp = $ws.Range("F715").Value2
$ws.Range("F715").Value = $ws.Range("F716").Value2
$ws.Range("F716").Value = $tmp
$tmp = $ws.Range("G715").Value2
$ws.Range("G715").Value = $ws.Range("G716").Value2
$ws.Range("G716").Value = $tmp

# Swap row 720 <-> row 721 for columns B,E,F,G
$tmp = $ws.Range("B720").Value2
$ws.Range("B720").Value = $ws.Range("B721").Value2
$ws.Range("B721").Value = $tmp
$tmp = $ws.Range("E720").Value2
$ws.Range("E720").Value = $ws.Range("E721").Value2
$ws.Range("E721").Value = $tmp
$tmp = $ws.Range("F720").Value2
$ws.Range("F720").Value = $ws.Range("F721").Value2
$ws.Range("F721").Value = $tmp
$tmp = $ws.Range("G720").Value2
$ws.Range("G720").Value = $ws.Range("G721").Value2
$ws.Range("G721").Value = $tmp
